$wb = $excel.ActiveWorkbook

# Fix the typo in the second sheet's name: "SNOWMED" -> "SNOMED"
$wsInclude = $wb.Worksheets.Item("Include from SNOWMED CT")
$wsInclude.Name = "Include from SNOMED CT"

# On the Metadata sheet: set "Experimental" value to "true" (as text, matching
# the original file's style) and refresh the "Date" value to the new timestamp.
$wsMeta = $wb.Worksheets.Item("Metadata")

# Use a leading apostrophe so the literal text "true" is stored (otherwise the
# host auto-converts the bareword to a real Boolean), then re-apply the
# original cell formatting on top so the style index is preserved.
$wsMeta.Cells.Item(7, 2).Value = "'true"
$wsMeta.Cells.Item(6, 2).Copy()
$wsMeta.Cells.Item(7, 2).PasteSpecial(-4122)

$wsMeta.Cells.Item(8, 2).Value = "2024-05-23T15:19:20+10:00"
